$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - DANIELA FLOREZ CAMARGO, period 2507
$ws.Range("C16").Value = "1047463262"
$ws.Range("D16").Value = "DANIELA FLOREZ CAMARGO"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 57520
$ws.Range("G16").Value = 1438000

# Row 17 - KATYA CAROLINA HAWKINS RAMIREZ, period 2507 (unchanged values, kept as-is)
$ws.Range("C17").Value = "1143164560"
$ws.Range("D17").Value = "KATYA CAROLINA HAWKINS RAMIREZ"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18 - DANIELA FLOREZ CAMARGO, period 2508
$ws.Range("C18").Value = "1047463262"
$ws.Range("D18").Value = "DANIELA FLOREZ CAMARGO"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 57520
$ws.Range("G18").Value = 1438000

# Row 19 - KATYA CAROLINA HAWKINS RAMIREZ, period 2508
$ws.Range("C19").Value = "1143164560"
$ws.Range("D19").Value = "KATYA CAROLINA HAWKINS RAMIREZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Totals
$ws.Range("E11").Value = 228920
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
